# Updates the cryptos list worksheet with the latest fetched values (prices & % volumes),
# including two coin rows that changed rank/order (TRON<->Polkadot, BitcoinSV<->Aptos).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.659.78"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "1.950.65"
$ws.Range("E3").Value = "  +1.42%  "
$ws.Range("D4").Value = "'0.9968"
$ws.Range("E4").Value = "  -0.38%  "
$ws.Range("D5").Value = "'246.70"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").Value = "'0.9971"
$ws.Range("E6").Value = "  -0.33%  "
$ws.Range("D7").Value = "'0.4842"
$ws.Range("E7").Value = "  +2.42%  "
$ws.Range("D8").Value = "'0.2918"
$ws.Range("E8").Value = "  +0.19%  "
$ws.Range("D9").Value = "'0.06836"
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("D10").Value = "'112.36"
$ws.Range("E10").Value = "  +6.36%  "
$ws.Range("D11").Value = "'19.65"
$ws.Range("E11").Value = "  +7.00%  "
$ws.Range("D12").Value = "1.935.94"
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.07583"
$ws.Range("E13").Value = "  -1.81%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'5.482"
$ws.Range("E14").Value = "  +2.93%  "
$ws.Range("D15").Value = "'0.6829"
$ws.Range("E15").Value = "  +1.88%  "
$ws.Range("D16").Value = "'301.74"
$ws.Range("E16").Value = "  +3.04%  "
$ws.Range("D17").Value = "30.615.37"
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("E18").Value = "  +1.53%  "
$ws.Range("D19").Value = "'0.000007698"
$ws.Range("E19").Value = "  +0.95%  "
$ws.Range("D20").Value = "'5.603"
$ws.Range("E20").Value = "  +0.74%  "
$ws.Range("D21").Value = "'0.9966"
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("D22").Value = "2.183.71"
$ws.Range("E22").Value = "  +0.49%  "
$ws.Range("D23").Value = "'0.9966"
$ws.Range("E23").Value = "  -0.43%  "
$ws.Range("D24").Value = "'6.531"
$ws.Range("E24").Value = "  +0.63%  "
$ws.Range("D25").Value = "'9.544"
$ws.Range("E25").Value = "  +0.27%  "
$ws.Range("D26").Value = "'167.88"
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("D27").Value = "'20.60"
$ws.Range("D28").Value = "'2.162"
$ws.Range("E28").Value = "  +1.45%  "
$ws.Range("D29").Value = "'0.1076"
$ws.Range("E29").Value = "  +0.52%  "
$ws.Range("D30").Value = "'1.441"
$ws.Range("E30").Value = "  +2.48%  "
$ws.Range("D31").Value = "'4.176"
$ws.Range("E31").Value = "  -0.50%  "
$ws.Range("D32").Value = "'4.103"
$ws.Range("E32").Value = "  +1.11%  "
$ws.Range("D33").Value = "'0.05004"
$ws.Range("E33").Value = "  -0.59%  "
$ws.Range("D34").Value = "'0.7450"
$ws.Range("E34").Value = "  +1.35%  "
$ws.Range("D35").Value = "'1.156"
$ws.Range("E35").Value = "  +0.97%  "
$ws.Range("D36").Value = "'0.02046"
$ws.Range("E36").Value = "  -0.74%  "
$ws.Range("D37").Value = "'2.708"
$ws.Range("E37").Value = "  -0.83%  "
$ws.Range("D38").Value = "'2.707"
$ws.Range("E38").Value = "  +0.84%  "
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("D40").Value = "'110.46"
$ws.Range("E40").Value = "  -1.12%  "
$ws.Range("D41").Value = "'0.4489"
$ws.Range("E41").Value = "  +1.06%  "
$ws.Range("D42").Value = "'0.8739"
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("D43").Value = "'5.876"
$ws.Range("E43").Value = "  -0.17%  "
$ws.Range("D44").Value = "'70.13"
$ws.Range("E44").Value = "  +3.17%  "
$ws.Range("D45").Value = "'1.000"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("B46").Value = "BitcoinSV"
$ws.Range("C46").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D46").Value = "'49.60"
$ws.Range("E46").Value = "  +3.59%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "'7.315"
$ws.Range("E47").Value = "  +0.25%  "
$ws.Range("D48").Value = "'9.363"
$ws.Range("E48").Value = "  -0.66%  "
$ws.Range("D49").Value = "'0.1241"
$ws.Range("E49").Value = "  -0.83%  "
$ws.Range("D50").Value = "'0.2550"
$ws.Range("E50").Value = "  +2.48%  "
$ws.Range("D51").Value = "'35.12"
$ws.Range("E51").Value = "  -0.17%  "
